$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 2, shifting the rest of the
# data (old rows 2-21) down to rows 5-24.
$ws.Rows("2:4").Insert()

# Inserting rows copies the formatting of the row above (the bold header),
# which the target data rows should not have - strip it back to "no style"
# to match the rest of the plain numeric data rows.
$ws.Range("A2:C4").ClearFormats()

# Fill in the 3 newly inserted rows (new rows 2-4) with the new data.
$ws.Range("A2").Value = -2.211326163574302
$ws.Range("B2").Value = 2.150324112876168
$ws.Range("C2").Value = 1.247732900374419

$ws.Range("A3").Value = -4.031598659201058
$ws.Range("B3").Value = 4.015500838530144
$ws.Range("C3").Value = -0.2885845313525067

$ws.Range("A4").Value = -6.420986591104682
$ws.Range("B4").Value = 1.20337057380033
$ws.Range("C4").Value = -1.732376253138717

# Append 7 brand-new rows (25-31) after the (shifted) existing data which
# now ends at row 24.
$ws.Range("A25").Value = 37.32190535987588
$ws.Range("B25").Value = -33.90012039951871
$ws.Range("C25").Value = 27.08640504549361

$ws.Range("A26").Value = -0.3525614285610033
$ws.Range("B26").Value = -16.87390499541134
$ws.Range("C26").Value = -5.992740513892805

$ws.Range("A27").Value = -27.24396839887731
$ws.Range("B27").Value = 3.699085454034913
$ws.Range("C27").Value = -26.36042895396987

$ws.Range("A28").Value = -11.78973249882937
$ws.Range("B28").Value = -39.63500136902936
$ws.Range("C28").Value = 31.93103018819374

$ws.Range("A29").Value = -3.788272580621008
$ws.Range("B29").Value = -24.42989606164616
$ws.Range("C29").Value = 38.63967183448689

$ws.Range("A30").Value = 10.84635844310596
$ws.Range("B30").Value = 4.736711256996685
$ws.Range("C30").Value = 13.0294989154327

$ws.Range("A31").Value = -15.62121211217174
$ws.Range("B31").Value = -12.64114606846652
$ws.Range("C31").Value = -11.19534088646288
